$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data for rows 2-7, columns A-T.
$rows = @{
    2 = @("FAPs","Pspn","Gfra4","ECs",3,1,3.042572666666667,9.127718,0.9827375016055572,0.9827375016055572,2,0.6666666666666666,0.2615776666666667,0.784733,0.09693893622669539,0.0969389362266954,0.7958690588104443,7.162821529294,0.09526552799572306,0.09526552799572308)
    3 = @("FAPs","Pspn","Gfra4","FAPs",3,1,3.042572666666667,9.127718,0.9827375016055572,0.9827375016055572,3,1,1.465947333333333,4.397842,0.5432702908799331,0.5432702908799331,4.460251287172889,40.142261584556,0.5338920883558698,0.5338920883558698)
    4 = @("FAPs","Pspn","Gfra4","sCs",3,1,3.042572666666667,9.127718,0.9827375016055572,0.9827375016055572,3,1,0.9708506666666668,2.912552,0.3597907728933716,0.3597907728933716,2.953883701815111,26.584953316336,0.3535798852539644,0.3535798852539644)
    5 = @("M2","Pspn","Gfra4","ECs",1,0.3333333333333333,0.053445,0.160335,0.01726249839444284,0.01726249839444284,2,0.6666666666666666,0.2615776666666667,0.784733,0.09693893622669539,0.0969389362266954,0.013980018395,0.125820165555,0.001673408230972326,0.001673408230972326)
    6 = @("M2","Pspn","Gfra4","FAPs",1,0.3333333333333333,0.053445,0.160335,0.01726249839444284,0.01726249839444284,3,1,1.465947333333333,4.397842,0.5432702908799331,0.5432702908799331,0.07834755522999999,0.70512799707,0.00937820252406334,0.00937820252406334)
    7 = @("M2","Pspn","Gfra4","sCs",1,0.3333333333333333,0.053445,0.160335,0.01726249839444284,0.01726249839444284,3,1,0.9708506666666668,2.912552,0.3597907728933716,0.3597907728933716,0.05188711388,0.4669840249200001,0.006210887639407176,0.006210887639407176)
}

foreach ($r in 2..7) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}
